$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

# Remove the "bankdeposits" row (row 9) entirely - shifts all following
# rows up by one, shrinks the table/used range, and drops the 3 shared
# strings ("bankdeposits", "Bank Deposits", "DPSACBW027SBOG") that were
# only referenced by that row.
$ws.Rows.Item(9).Delete()

$ws.Range("D11").Select()
